$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")

$r = $ws.Range("A39")
$r.Style = "Normal"
